$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# Add header "Correction" in N1, matching the style of existing headers (e.g. M1)
$ws.Range("N1").Value = "Correction"
$ws.Range("N1").Style = $ws.Range("M1").Style

# Fill M2:M12 with "nan" where currently empty, and add blank text cells in N2:N12
for ($r = 2; $r -le 12; $r++) {
    $mCell = $ws.Cells.Item($r, 13)  # column M = 13
    if ([string]::IsNullOrEmpty($mCell.Value)) {
        $mCell.Value = "nan"
    }
    $nCell = $ws.Cells.Item($r, 14)  # column N = 14
    $nCell.Value = ""
}
